# Auto-generated edit script: updates cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $oldStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $oldStyle
}

Set-TextValue $ws.Range("D2") '41.375.55'
Set-TextValue $ws.Range("E2") '  -1.05%  '
Set-TextValue $ws.Range("D3") '2.191.26'
Set-TextValue $ws.Range("E3") '  -0.96%  '
Set-TextValue $ws.Range("E4") '  +0.00%  '
Set-TextValue $ws.Range("D5") '255.14'
Set-TextValue $ws.Range("E5") '  +3.65%  '
Set-TextValue $ws.Range("E6") '  -0.46%  '
Set-TextValue $ws.Range("D7") '67.98'
Set-TextValue $ws.Range("E7") '  -3.80%  '
Set-TextValue $ws.Range("E8") '  -0.03%  '
Set-TextValue $ws.Range("E9") '  +5.06%  '
Set-TextValue $ws.Range("D10") '37.89'
Set-TextValue $ws.Range("E10") '  +2.53%  '
Set-TextValue $ws.Range("D11") '58.13'
Set-TextValue $ws.Range("E11") '  +0.21%  '
Set-TextValue $ws.Range("E12") '  -1.17%  '
Set-TextValue $ws.Range("E13") '  +5.16%  '
Set-TextValue $ws.Range("E14") '  -1.08%  '
Set-TextValue $ws.Range("D15") '2.518.04'
Set-TextValue $ws.Range("E15") '  -0.94%  '
Set-TextValue $ws.Range("D16") '0.871'
Set-TextValue $ws.Range("E16") '  +2.96%  '
Set-TextValue $ws.Range("D17") '14.52'
Set-TextValue $ws.Range("E17") '  -2.38%  '
Set-TextValue $ws.Range("D18") '2.214.56'
Set-TextValue $ws.Range("E18") '  +0.05%  '
Set-TextValue $ws.Range("D19") '41.251.84'
Set-TextValue $ws.Range("E19") '  -1.16%  '
Set-TextValue $ws.Range("D20") '0.0₃0953'
Set-TextValue $ws.Range("E20") '  -0.65%  '
Set-TextValue $ws.Range("D21") '6.25'
Set-TextValue $ws.Range("E21") '  +2.19%  '
Set-TextValue $ws.Range("D22") '71.97'
Set-TextValue $ws.Range("E22") '  -2.21%  '
Set-TextValue $ws.Range("D23") '232.70'
Set-TextValue $ws.Range("E23") '  -1.14%  '
Set-TextValue $ws.Range("D24") '2.08'
Set-TextValue $ws.Range("E24") '  +1.38%  '
Set-TextValue $ws.Range("D25") '12.00'
Set-TextValue $ws.Range("E25") '  +20.94%  '
Set-TextValue $ws.Range("D26") '3.81'
Set-TextValue $ws.Range("E26") '  +5.65%  '
Set-TextValue $ws.Range("E27") '  +0.01%  '
Set-TextValue $ws.Range("D28") '2.52'
Set-TextValue $ws.Range("E28") '  +2.68%  '
Set-TextValue $ws.Range("E29") '  -2.02%  '
Set-TextValue $ws.Range("D30") '169.18'
Set-TextValue $ws.Range("E30") '  -0.58%  '
Set-TextValue $ws.Range("D31") '20.62'
Set-TextValue $ws.Range("E31") '  +0.70%  '
Set-TextValue $ws.Range("E32") '  +0.47%  '
Set-TextValue $ws.Range("E33") '  -2.44%  '
Set-TextValue $ws.Range("E34") '  +6.40%  '
Set-TextValue $ws.Range("D35") '0.0728'
Set-TextValue $ws.Range("E35") '  +1.44%  '
Set-TextValue $ws.Range("D36") '4.61'
Set-TextValue $ws.Range("E36") '  -0.75%  '
Set-TextValue $ws.Range("D37") '25.17'
Set-TextValue $ws.Range("E37") '  +8.75%  '
Set-TextValue $ws.Range("D38") '4.00'
Set-TextValue $ws.Range("E38") '  +3.87%  '
Set-TextValue $ws.Range("E39") '  +8.92%  '
Set-TextValue $ws.Range("E40") '  -2.07%  '
Set-TextValue $ws.Range("D41") '5.75'
Set-TextValue $ws.Range("E41") '  -2.08%  '
Set-TextValue $ws.Range("D42") '12.26'
Set-TextValue $ws.Range("E42") '  +18.17%  '
Set-TextValue $ws.Range("D43") '64.36'
Set-TextValue $ws.Range("E43") '  -2.39%  '
Set-TextValue $ws.Range("E44") '  +7.47%  '
Set-TextValue $ws.Range("D45") '4.88'
Set-TextValue $ws.Range("E45") '  -0.19%  '
Set-TextValue $ws.Range("B46") 'FraxShare'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D46") '8.63'
Set-TextValue $ws.Range("E46") '  -3.91%  '
Set-TextValue $ws.Range("B47") 'Cronos'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D47") '0.101'
Set-TextValue $ws.Range("E47") '  +1.66%  '
Set-TextValue $ws.Range("D48") '1.00'
Set-TextValue $ws.Range("E48") '  -0.01%  '
Set-TextValue $ws.Range("D49") '1.15'
Set-TextValue $ws.Range("E49") '  +4.06%  '
Set-TextValue $ws.Range("E50") '  -0.66%  '
Set-TextValue $ws.Range("E51") '  +0.27%  '
